# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" sheet, which mirror the same event rows.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F10" = 588
    "F12" = 313
    "F39" = 3760
    "F40" = 437
    "F45" = 77
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
